$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 14 (shifting C14..Y1 down by one row)
$ws.Rows.Item(14).Insert()

# Fill in the new row 14 with C13 data
$ws.Cells.Item(14, 1).Value = "C13"
$ws.Cells.Item(14, 2).Value = 127.7
$ws.Cells.Item(14, 3).Value = -112.3
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = "top"

# Copy style from row 13/15 (B:D columns use style s=1, numFmt 100) to the new row
$ws.Range("B13:D13").Copy()
$ws.Range("B14:D14").PasteSpecial(-4122) | Out-Null

# Re-apply the numeric formatting over the whole B:D data range so every
# cell (including the row shifted in by the insert) keeps the "0.000000"
# number format consistently.
$ws.Range("B1:D55").NumberFormat = "0.000000"

# Set A1 selection
$ws.Range("A1:E1").Select()
